$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header B1 "value" -> "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# New date (column A) / value (column B) series. The series starts/ends
# with a date-only row (no release value yet available for that year).
$dates  = @(38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657, 46022)
$values = @($null, 3.083829270092076, 2.599189965723969, 0.9885759521669257, -5.06641082168553, 2.964960767277169, 3.148381490172691, 0.8410377950035519, 0.4885704456499607, 1.581181299492873, 1.369290370184317, 1.652444083200288, 2.175598592631678, 1.570309077096765, 0.4185134408791091, -6.428698128894917, 1.195045614048973, 1.373420289353544, -0.2994265979745614, 0.1140152762751701, $null)

$lastRow = 1 + $dates.Length

# Extend the date formatting (copied from A2, which already carries the
# correct number format/border/font) down through the new rows before
# writing values, so the new cells pick up the same style.
$ws.Range("A2").Copy()
$ws.Range("A3:A$lastRow").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

Write-Host "done"
